# smaller fixes to resume
#
# 1. Update the "Sportteri.fi" experience entry end-date on the Kokemus sheet
#    from 01/2014 - 06/2017 to 01/2014 - 06/2018.
# 2. Fix the "Datan kasittely" skill description on the Taidot sheet
#    (typo swap "tai" -> "js").
# 3. Update view/selection state: Koulutus and Kokemus sheets get new
#    (non-active) cell selections, while Taidot becomes the active sheet
#    with its own new selection.

$wb = $excel.ActiveWorkbook

# --- Content fixes -----------------------------------------------------

$wsKokemus = $wb.Worksheets.Item("Kokemus")
$wsKokemus.Range("D4").Value2 = "01/2014 - 06/2018"

$wsTaidot = $wb.Worksheets.Item("Taidot")
$wsTaidot.Range("B6").Value2 = "Olen työskennellyt suurten tietokantojen ja tietolähteiden kanssa. Esimerkiksi tietokannan siirto, web-karttatietojen käsittely (GeoJSON KML, SHP) js kuvankäsittelyn automatisointi."

# --- View / selection fixes --------------------------------------------

$wsKoulutus = $wb.Worksheets.Item("Koulutus")
$wsKoulutus.Select()
$wsKoulutus.Range("G4").Select()

$wsKokemus.Select()
$wsKokemus.Range("F5").Select()

$wsTaidot.Select()
$wsTaidot.Range("J5").Select()
